$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Execute" (Y) markers from the Quarterly Estimated Tax and
# New Tax Return Amount Due rows, leaving only the Extension Payments row
# flagged for execution (RAD Extension Payment Type run).
$ws.Range("C2").Clear()
$ws.Range("C4").Clear()

# Record the new execution timestamp for the Extension Payments row.
$ws.Range("B3").Value = "Wed Mar 20 23:05:59 EDT 2024"
